# CEDS Data and exogenous assumptions updates
# Insert a new sector row ("1A1bc_Other-feedstocks") into the "Sectors" sheet,
# right after "1A1bc_Other-transformation" (row 5), pushing every row below
# it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new blank row at row 6 (shifts rows 6..60 down to 7..61,
# carrying along formatting/styles as Excel normally does).
$ws.Rows.Item(6).Insert()

# Populate the new row with the new sector's data.
$ws.Cells.Item(6, 1).Value = "1A1bc_Other-feedstocks"
$ws.Cells.Item(6, 2).Value = "Energy_Combustion"
$ws.Cells.Item(6, 3).Value = "kt"
$ws.Cells.Item(6, 4).Value = "NC"

# Update the selection to reflect the new active cell (A6) as in the saved file.
$ws.Range("A6").Select()
